$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) and Volume(1h) (column E) updates for the refreshed crypto snapshot.
# Numeric-looking Price values must be written as text (matching the source data, which
# stores every Price cell as a plain string) by forcing a Text number format before the
# assignment - otherwise Excel COM auto-converts them to numbers and trailing zeros
# (e.g. "86.00") would be silently dropped.

$ws.Range("D2").Value = "44.015.51"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "2.237.09"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.67"
$ws.Range("E5").Value = "  -3.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.05"
$ws.Range("E6").Value = "  -5.72%  "
$ws.Range("E7").Value = "  -1.06%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  -5.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.88"
$ws.Range("E10").Value = "  -5.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0807"
$ws.Range("E11").Value = "  -3.13%  "
$ws.Range("E12").Value = "  -4.55%  "
$ws.Range("E13").Value = "  -1.56%  "
$ws.Range("D14").Value = "2.578.01"
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").Value = "2.237.19"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("E16").Value = "  -3.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.58"
$ws.Range("E17").Value = "  -6.25%  "
$ws.Range("D18").Value = "43.913.14"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "0.0₃0960"
$ws.Range("E19").Value = "  -2.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.20"
$ws.Range("E20").Value = "  -8.58%  "
$ws.Range("E21").Value = "  -3.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.01"
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.26"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("E24").Value = "  -5.87%  "
$ws.Range("E25").Value = "  -5.14%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -6.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.63"
$ws.Range("E28").Value = "  -2.42%  "
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.95"
$ws.Range("E30").Value = "  -2.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.88"
$ws.Range("E31").Value = "  -1.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "154.20"
$ws.Range("E32").Value = "  -4.32%  "
$ws.Range("E33").Value = "  -4.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.21"
$ws.Range("E34").Value = "  +3.90%  "
$ws.Range("E35").Value = "  -3.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.118"
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.106"
$ws.Range("E37").Value = "  -5.56%  "
$ws.Range("E38").Value = "  -9.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.10"
$ws.Range("E39").Value = "  -7.97%  "
$ws.Range("E40").Value = "  -8.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.81"
$ws.Range("E41").Value = "  -7.78%  "
$ws.Range("E42").Value = "  -4.37%  "
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").Value = "1.737.44"
$ws.Range("E44").Value = "  -1.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "86.00"
$ws.Range("E45").Value = "  +6.50%  "
$ws.Range("E46").Value = "  -4.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "99.86"
$ws.Range("E47").Value = "  -3.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.91"
$ws.Range("E48").Value = "  -5.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.10"
$ws.Range("E49").Value = "  -2.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "69.09"
$ws.Range("E50").Value = "  -7.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.22"
$ws.Range("E51").Value = "  -6.17%  "
